# Natmi following Dr Hou advice
# Update LR-pairs (Wnt2-Fzd8) sheet: refresh computed statistics for the
# existing target clusters, relabel the former "Neutro"/"sCs" rows as the
# new "M1"/"M2" macrophage-polarization clusters, and append fresh rows
# for the actual "Neutro" and "sCs" clusters.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2 (Target cluster: ECs) ---
$ws.Cells.Item(2, 5).Value  = 3
$ws.Cells.Item(2, 6).Value  = 1
$ws.Cells.Item(2, 7).Value  = 0.574538
$ws.Cells.Item(2, 8).Value  = 1.723614
$ws.Cells.Item(2, 9).Value  = 1
$ws.Cells.Item(2, 10).Value = 1
$ws.Cells.Item(2, 11).Value = 2
$ws.Cells.Item(2, 12).Value = 1
$ws.Cells.Item(2, 13).Value = 2.569008
$ws.Cells.Item(2, 14).Value = 5.138016
$ws.Cells.Item(2, 15).Value = 0.2577238367773512
$ws.Cells.Item(2, 16).Value = 0.2133397199412101
$ws.Cells.Item(2, 17).Value = 1.475992718304
$ws.Cells.Item(2, 18).Value = 8.855956309824
$ws.Cells.Item(2, 19).Value = 0.2577238367773512
$ws.Cells.Item(2, 20).Value = 0.2133397199412101

# --- Row 3 (Target cluster: FAPs) ---
$ws.Cells.Item(3, 5).Value  = 3
$ws.Cells.Item(3, 6).Value  = 1
$ws.Cells.Item(3, 7).Value  = 0.574538
$ws.Cells.Item(3, 8).Value  = 1.723614
$ws.Cells.Item(3, 9).Value  = 1
$ws.Cells.Item(3, 10).Value = 1
$ws.Cells.Item(3, 11).Value = 3
$ws.Cells.Item(3, 12).Value = 1
$ws.Cells.Item(3, 13).Value = 2.043481333333333
$ws.Cells.Item(3, 14).Value = 6.130444
$ws.Cells.Item(3, 15).Value = 0.2050028063787906
$ws.Cells.Item(3, 16).Value = 0.2545471259870097
$ws.Cells.Item(3, 17).Value = 1.174057678290667
$ws.Cells.Item(3, 18).Value = 10.566519104616
$ws.Cells.Item(3, 19).Value = 0.2050028063787906
$ws.Cells.Item(3, 20).Value = 0.2545471259870097

# --- Row 4 (Target cluster relabeled Neutro -> M1) ---
$ws.Cells.Item(4, 4).Value  = "M1"
$ws.Cells.Item(4, 5).Value  = 3
$ws.Cells.Item(4, 6).Value  = 1
$ws.Cells.Item(4, 7).Value  = 0.574538
$ws.Cells.Item(4, 8).Value  = 1.723614
$ws.Cells.Item(4, 9).Value  = 1
$ws.Cells.Item(4, 10).Value = 1
$ws.Cells.Item(4, 11).Value = 2
$ws.Cells.Item(4, 12).Value = 0.6666666666666666
$ws.Cells.Item(4, 13).Value = 0.01206
$ws.Cells.Item(4, 14).Value = 0.03618
$ws.Cells.Item(4, 15).Value = 0.001209863679496076
$ws.Cells.Item(4, 16).Value = 0.001502259056311421
$ws.Cells.Item(4, 17).Value = 0.00692892828
$ws.Cells.Item(4, 18).Value = 0.06236035451999999
$ws.Cells.Item(4, 19).Value = 0.001209863679496076
$ws.Cells.Item(4, 20).Value = 0.001502259056311421

# --- Row 5 (Target cluster relabeled sCs -> M2) ---
$ws.Cells.Item(5, 4).Value  = "M2"
$ws.Cells.Item(5, 5).Value  = 3
$ws.Cells.Item(5, 6).Value  = 1
$ws.Cells.Item(5, 7).Value  = 0.574538
$ws.Cells.Item(5, 8).Value  = 1.723614
$ws.Cells.Item(5, 9).Value  = 1
$ws.Cells.Item(5, 10).Value = 1
$ws.Cells.Item(5, 11).Value = 2
$ws.Cells.Item(5, 12).Value = 0.6666666666666666
$ws.Cells.Item(5, 13).Value = 0.03685233333333333
$ws.Cells.Item(5, 14).Value = 0.110557
$ws.Cells.Item(5, 15).Value = 0.003697039768215801
$ws.Cells.Item(5, 16).Value = 0.004590526658060303
$ws.Cells.Item(5, 17).Value = 0.02117306588866667
$ws.Cells.Item(5, 18).Value = 0.190557592998
$ws.Cells.Item(5, 19).Value = 0.003697039768215801
$ws.Cells.Item(5, 20).Value = 0.004590526658060303

# --- Row 6 (new row, Target cluster: Neutro) ---
$ws.Cells.Item(6, 1).Value  = "FAPs"
$ws.Cells.Item(6, 2).Value  = "Wnt2"
$ws.Cells.Item(6, 3).Value  = "Fzd8"
$ws.Cells.Item(6, 4).Value  = "Neutro"
$ws.Cells.Item(6, 5).Value  = 3
$ws.Cells.Item(6, 6).Value  = 1
$ws.Cells.Item(6, 7).Value  = 0.574538
$ws.Cells.Item(6, 8).Value  = 1.723614
$ws.Cells.Item(6, 9).Value  = 1
$ws.Cells.Item(6, 10).Value = 1
$ws.Cells.Item(6, 11).Value = 3
$ws.Cells.Item(6, 12).Value = 1
$ws.Cells.Item(6, 13).Value = 2.055205
$ws.Cells.Item(6, 14).Value = 6.165615
$ws.Cells.Item(6, 15).Value = 0.206178928973361
$ws.Cells.Item(6, 16).Value = 0.2560074895378535
$ws.Cells.Item(6, 17).Value = 1.18079337029
$ws.Cells.Item(6, 18).Value = 10.62714033261
$ws.Cells.Item(6, 19).Value = 0.206178928973361
$ws.Cells.Item(6, 20).Value = 0.2560074895378535

# --- Row 7 (new row, Target cluster: sCs) ---
$ws.Cells.Item(7, 1).Value  = "FAPs"
$ws.Cells.Item(7, 2).Value  = "Wnt2"
$ws.Cells.Item(7, 3).Value  = "Fzd8"
$ws.Cells.Item(7, 4).Value  = "sCs"
$ws.Cells.Item(7, 5).Value  = 3
$ws.Cells.Item(7, 6).Value  = 1
$ws.Cells.Item(7, 7).Value  = 0.574538
$ws.Cells.Item(7, 8).Value  = 1.723614
$ws.Cells.Item(7, 9).Value  = 1
$ws.Cells.Item(7, 10).Value = 1
$ws.Cells.Item(7, 11).Value = 2
$ws.Cells.Item(7, 12).Value = 1
$ws.Cells.Item(7, 13).Value = 3.2514585
$ws.Cells.Item(7, 14).Value = 6.502917
$ws.Cells.Item(7, 15).Value = 0.3261875244227855
$ws.Cells.Item(7, 16).Value = 0.2700128788195549
$ws.Cells.Item(7, 17).Value = 1.868086463673
$ws.Cells.Item(7, 18).Value = 11.208518782038
$ws.Cells.Item(7, 19).Value = 0.3261875244227855
$ws.Cells.Item(7, 20).Value = 0.2700128788195549
